$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.971.43'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '2.553.93'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.87'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.26'
$ws.Range('E6').Value = '  +6.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.579'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.21'
$ws.Range('E10').Value = '  +3.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0827'
$ws.Range('E11').Value = '  +2.62%  '
$ws.Range('E12').Value = '  +4.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.64'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('D14').Value = '2.946.26'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = '2.548.13'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.06'
$ws.Range('E16').Value = '  +6.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.876'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').Value = '43.008.62'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.86'
$ws.Range('E19').Value = '  +5.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.60'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.01'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '254.59'
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '28.10'
$ws.Range('E26').Value = '  -3.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.22'
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.90'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.09'
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '159.06'
$ws.Range('E32').Value = '  +3.24%  '
$ws.Range('B33').Value = 'Celestia'
$ws.Range('C33').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.55'
$ws.Range('E33').Value = '  +15.00%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.16'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.75'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0804'
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.31'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.75'
$ws.Range('E39').Value = '  +11.14%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.43'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.08'
$ws.Range('E42').Value = '  +30.96%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('D45').Value = '2.084.69'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.86'
$ws.Range('E47').Value = '  +2.32%  '
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('D49').Value = '2.803.45'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.86'
$ws.Range('E50').Value = '  +8.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.37'
$ws.Range('E51').Value = '  -1.16%  '
